# Routing-Test-Scripts.xlsx update
# - Swap the Panorama-2 public IP (B6) for a freshly reallocated AWS address
#   and shrink its font to match the other "pub" IP cells (B3:B5).
# - Mark several ping checks that came back negative with "NO" in column D.
# - Add a purple "how to use this spreadsheet" instructions banner in F1:Q8.
# - Leave the active selection on A8 (where the new instructions start).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$purple = 10498160   # RGB(112,48,160) -> FF7030A0
$white  = 16777215   # RGB(255,255,255)

# --- 1) Panorama-2 public IP changed ---------------------------------------
$ws.Range("B6").Value = "52.11.212.64"
$ws.Range("B6").Font.Name = "Times New Roman"
$ws.Range("B6").Font.Size = 5
$ws.Range("B6").Font.Color = 0

# --- 2) Flag the failed ping checks with "NO" in column D ------------------
$ws.Range("D39").Value = "NO"
$ws.Range("D40").Value = "NO"
$ws.Range("D44").Value = "NO"
$ws.Range("D45").Value = "NO"
$ws.Range("D52").Value = "NO"
$ws.Range("D53").Value = "NO"
$ws.Range("D57").Value = "NO"
$ws.Range("D58").Value = "NO"

# --- 3) New "how to use this spreadsheet" banner, F1:Q8 ---------------------
$banner = $ws.Range("F1:Q8")
$banner.Interior.Color = $purple
$banner.Font.Name = "Calibri"
$banner.Font.Size = 11
$banner.Font.Color = $white
$banner.Font.Bold = $false
$banner.Font.Underline = $false

$ws.Range("F1").Value = "How to use this spreadsheet to check all of your AWS routes:"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F1").Font.Underline = $true

$ws.Range("F2").Value = "1. Run all of the terraform build scripts and the RT-Associations.sh bash script"
$ws.Range("F2").Font.Bold = $true

$ws.Range("F3").Value = "2. Get all of the public IPs from AWS and paste them into the cells B3 - C6"
$ws.Range("F3").Font.Bold = $true

$ws.Range("F4").Value = "3. Download the new keypair file created by the terraform scripts and paste the name into B2."
$ws.Range("F4").Font.Bold = $true

$ws.Range("F5").Value = "4. Upload the keypair file into your ssh system to test from. A sample file name is :terraform-key-pair.1e5a.pem"
$ws.Range("F5").Font.Bold = $true

$ws.Range("F6").Value = "5. Run the attribute changes as in the sss-attributes.bat script (for windows), in linux use chmod"
$ws.Range("F6").Font.Bold = $true

$ws.Range("F7").Value = "6. Copy the ssh strings in blue, one at a time from your test system. Run the ssh, connect to the AWS system. "
$ws.Range("F7").Font.Bold = $true

$ws.Range("F8").Value = "7. Paste in all of the ping commands with IP addresses and check the results from your test host. Expected results are in column D. You can record actuals in column E if you chose. "
$ws.Range("F8").Font.Bold = $true

# --- 4) Leave the selection where the author left it ------------------------
$ws.Range("A8").Select()
